$d = $word.ActiveDocument

$pairs = @(
    @("2025-02-01 Saturday", "2025-02-02 Sunday"),
    @("793×5=3965", "491×4=1964"),
    @("849×7=5943", "192×6=1152"),
    @("285×4=1140", "770×8=6160"),
    @("817×2=1634", "184×6=1104"),
    @("816×3=2448", "945×8=7560"),
    @("905×6=5430", "296×2=592"),
    @("642×3=1926", "906×4=3624"),
    @("619×9=5571", "497×4=1988"),
    @("255×6=1530", "243×4=972"),
    @("550×3=1650", "742×6=4452"),
    @("719×6=4314", "111×2=222"),
    @("108×5=540", "439×6=2634"),
    @("749×7=5243", "961×8=7688"),
    @("527×3=1581", "461×2=922"),
    @("476×8=3808", "329×9=2961"),
    @("626×9=5634", "691×6=4146"),
    @("877×8=7016", "204×3=612"),
    @("395×4=1580", "848×7=5936"),
    @("809×5=4045", "332×4=1328"),
    @("690×2=1380", "333×8=2664"),
    @("672×4=2688", "216×2=432"),
    @("888×8=7104", "419×5=2095"),
    @("797×7=5579", "417×7=2919"),
    @("614×9=5526", "991×4=3964"),
    @("555×5=2775", "905×9=8145")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
